# Regenerate the handback status report: refreshed timestamps (report was
# re-run ~52s later) and one file's status flipped from "ht" (human
# translation) to "mt" (machine translation).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the two
# rows that shared the old timestamp.
$wsOverview.Range("G2").Value = "2016-08-16 16:13:34"
$wsOverview.Range("G5").Value = "2016-08-16 16:13:34"

# zh-cn sheet: Priority ("ht" -> "mt"), Correspond Handoff Datetime (H),
# Correspond Handback DateTime (K).
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-08-16 16:13:30"
$wsZhCn.Range("H5").Value = "2016-08-16 16:13:30"
$wsZhCn.Range("K2").Value = "2016-08-16 16:13:47"
$wsZhCn.Range("K5").Value = "2016-08-16 16:13:47"

# de-de sheet: Priority ("ht" -> "mt"), Correspond Handoff Datetime (H),
# Correspond Handback DateTime (K).
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-08-16 16:13:34"
$wsDeDe.Range("H5").Value = "2016-08-16 16:13:34"
$wsDeDe.Range("K2").Value = "2016-08-16 16:13:53"
$wsDeDe.Range("K5").Value = "2016-08-16 16:13:53"
